$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("62:62").Insert()
$ws.Range("A62").Value = 2
$ws.Range("B62").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 45134
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100112026
$ws.Range("G62").Value = "Haba"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 800
$ws.Range("K62").Value = 8500
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = 8750
$ws.Range("N62").Value = "$/saco 25 kilos"
$ws.Range("O62").Value = "Provincia de Limarí"
$ws.Range("P62").Value = 350
$ws.Range("Q62").Value = 25
$ws.Range("R62").Value = "Hortaliza"
